$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 40
$ws.Range("J11").Value = 0

# Row 12 updates (new data for day eleven)
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 20
$ws.Range("I12").Value = 60
$ws.Range("J12").Value = 0

# Update selection to reflect the active cell after edits
$ws.Range("C12").Select()
